$wb = $excel.ActiveWorkbook

# Update the BBNPPTY sheet values: set columns B, C, D (years 2021-2023) to 1
# for all fuel-type rows (2 through 25) to reflect banning new power plants
# in those years.
$ws = $wb.Worksheets.Item("BBNPPTY")
$ws.Range("B2:D25").Value = 1

# Make BBNPPTY the active sheet (and its selection/scroll position),
# and clear the tab-selected flag from the About sheet. Select on the
# About sheet first so that the final active sheet/tab ends up being
# BBNPPTY.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Select()
$wsAbout.Range("B18").Select()

$ws.Select()
$ws.Range("B2:D25").Select()
